# Automatic update of files.
# Bump the "Förändrad" (Changed) date column C for every data row
# (rows 2-106) from 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C106").Value = 45182
